# Apply the KHL injuries snapshot refresh:
#  - "snapshot" sheet: the player whose injury record returned (АДМ / Шепелев
#    Александр, old row 7) is removed from the open-injuries snapshot, and
#    every remaining row's "scraped_at" timestamp is refreshed to reflect the
#    new scrape run.
#  - "returned" sheet: the previous "returned" log (НХК/Попугаев and
#    ЦСК/Саморуков) is replaced by the single new return event
#    (АДМ/Шепелев Александр).
#  - "new_injured" sheet: unchanged (header only).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "snapshot"
# ---------------------------------------------------------------------------
$snap = $wb.Worksheets.Item("snapshot")

# Remove the row for АДМ / Шепелев Александр (row 7); every row below shifts
# up by one automatically, which matches the diff exactly.
$snap.Rows.Item(7).Delete()

# Refresh the "scraped_at" column (K) for every remaining data row (2-40)
# with the new scrape timestamps.
$scrapedAt = @{
    2  = "2025-11-19T03:06:26.519705+00:00"
    3  = "2025-11-19T03:06:26.519745+00:00"
    4  = "2025-11-19T03:06:26.519766+00:00"
    5  = "2025-11-19T03:06:31.003614+00:00"
    6  = "2025-11-19T03:06:31.003648+00:00"
    7  = "2025-11-19T03:06:36.636032+00:00"
    8  = "2025-11-19T03:06:42.103802+00:00"
    9  = "2025-11-19T03:06:47.664169+00:00"
    10 = "2025-11-19T03:06:47.664204+00:00"
    11 = "2025-11-19T03:06:58.676781+00:00"
    12 = "2025-11-19T03:07:04.141876+00:00"
    13 = "2025-11-19T03:07:09.208161+00:00"
    14 = "2025-11-19T03:07:09.208189+00:00"
    15 = "2025-11-19T03:07:09.208208+00:00"
    16 = "2025-11-19T03:07:14.311102+00:00"
    17 = "2025-11-19T03:07:19.862461+00:00"
    18 = "2025-11-19T03:07:25.555802+00:00"
    19 = "2025-11-19T03:07:30.748719+00:00"
    20 = "2025-11-19T03:07:30.748748+00:00"
    21 = "2025-11-19T03:07:30.748767+00:00"
    22 = "2025-11-19T03:07:30.748783+00:00"
    23 = "2025-11-19T03:07:36.331556+00:00"
    24 = "2025-11-19T03:07:36.331585+00:00"
    25 = "2025-11-19T03:07:41.532617+00:00"
    26 = "2025-11-19T03:07:41.532648+00:00"
    27 = "2025-11-19T03:07:41.532671+00:00"
    28 = "2025-11-19T03:07:47.098378+00:00"
    29 = "2025-11-19T03:07:47.098404+00:00"
    30 = "2025-11-19T03:07:52.165812+00:00"
    31 = "2025-11-19T03:07:52.165847+00:00"
    32 = "2025-11-19T03:07:52.165865+00:00"
    33 = "2025-11-19T03:07:52.165881+00:00"
    34 = "2025-11-19T03:07:52.165895+00:00"
    35 = "2025-11-19T03:07:57.271425+00:00"
    36 = "2025-11-19T03:07:57.271456+00:00"
    37 = "2025-11-19T03:08:08.347526+00:00"
    38 = "2025-11-19T03:08:08.347560+00:00"
    39 = "2025-11-19T03:08:08.347583+00:00"
    40 = "2025-11-19T03:08:13.846901+00:00"
}

foreach ($row in $scrapedAt.Keys) {
    $snap.Cells.Item($row, 11).Value = $scrapedAt[$row]
}

# ---------------------------------------------------------------------------
# Sheet "returned"
# ---------------------------------------------------------------------------
$ret = $wb.Worksheets.Item("returned")

# Drop the old second data row (ЦСК / Саморуков Дмитрий); the remaining
# row (НХК / Попугаев) will be overwritten below with the new return event.
$ret.Rows.Item(3).Delete()

$ret.Cells.Item(2, 1).Value = "АДМ"
$ret.Cells.Item(2, 2).Value = "Адмирал"
$ret.Cells.Item(2, 3).Value = "Шепелев Александр"
$ret.Cells.Item(2, 4).Value = "1369_АДМ_шепелевалександр"
$ret.Cells.Item(2, 5).Value = "RETURN"
$ret.Cells.Item(2, 6).Value = "2025-11-19T11:08:14.354550+08:00"

# Force the "changed_day" cell to stay plain text (otherwise Excel would
# auto-convert the yyyy-mm-dd looking string into a date serial number).
$ret.Cells.Item(2, 7).NumberFormat = "@"
$ret.Cells.Item(2, 7).Value = "2025-11-19"
